# EDA_Preprocess.pptx - "Add files via upload" edit
#
# 1. The "Today" date placeholder (a datetime2 field cached by PowerPoint)
#    on the slide master and on every slide layout advanced from
#    "Wednesday, October 23, 2024" to "Thursday, October 24, 2024".
# 2. On slide 2, the bullet describing the dataset was updated to mention
#    cancellation details in addition to booking details.

$p = $ppt.ActivePresentation

$oldDate = "Wednesday, October 23, 2024"
$newDate = "Thursday, October 24, 2024"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $phType = $null
        try { $phType = $shp.PlaceholderFormat.Type } catch { $phType = $null }
        if ($phType -eq $ppPlaceholderDate) {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master's own date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 2: mention cancellation details as well as booking details.
$s2 = $p.Slides.Item(2)
$contentShape = $s2.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$paragraphs = $tr.Paragraphs()
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like "The dataset: booking details*") {
        $para.Text = "The dataset: booking & cancellation details for city and resort hotels, featuring variables like booking dates, length of stay, guest count, # of special requests, …, with all personal information removed"
        break
    }
}
